# Add displacement and reaction calculation to structure analysis
# - Adds "F", "Desplazamientos" and "Reacciones" sheets after "Global"
# - Rebuilds the local stiffness matrices on "Barra 1 (K1)" / "Barra 2 (K2)"
#   and the assembled matrix on "Global" for correct 3D frame behavior

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fix the local 3D-frame stiffness matrices for Barra 1 (K1) and
#    Barra 2 (K2), and the assembled Global stiffness matrix.
# ---------------------------------------------------------------------
$k12x12 = @(
    @(4032,0,0,0,10080,0,-4032,0,0,0,10080,0),
    @(0,2016,0,-5040,0,0,0,-2016,0,-5040,0,0),
    @(0,0,42000,0,0,0,0,0,-42000,0,0,0),
    @(0,-5040,0,16800,0,0,0,5040,0,8400,0,0),
    @(10080,0,0,0,33600,0,-10080,0,0,0,16800,0),
    @(0,0,0,0,0,800,0,0,0,0,0,-800),
    @(-4032,0,0,0,-10080,0,4032,0,0,0,-10080,0),
    @(0,-2016,0,5040,0,0,0,2016,0,5040,0,0),
    @(0,0,-42000,0,0,0,0,0,42000,0,0,0),
    @(0,-5040,0,8400,0,0,0,5040,0,16800,0,0),
    @(10080,0,0,0,16800,0,-10080,0,0,0,33600,0),
    @(0,0,0,0,0,-800,0,0,0,0,0,800)
  )

$kGlobal18x18 = @(
    @(4032,0,0,0,10080,0,-4032,0,0,0,10080,0,0,0,0,0,0,0),
    @(0,2016,0,-5040,0,0,0,-2016,0,-5040,0,0,0,0,0,0,0,0),
    @(0,0,42000,0,0,0,0,0,-42000,0,0,0,0,0,0,0,0,0),
    @(0,-5040,0,16800,0,0,0,5040,0,8400,0,0,0,0,0,0,0,0),
    @(10080,0,0,0,33600,0,-10080,0,0,0,16800,0,0,0,0,0,0,0),
    @(0,0,0,0,0,800,0,0,0,0,0,-800,0,0,0,0,0,0),
    @(-4032,0,0,0,-10080,0,8064,0,0,0,0,0,-4032,0,0,0,10080,0),
    @(0,-2016,0,5040,0,0,0,4032,0,0,0,0,0,-2016,0,-5040,0,0),
    @(0,0,-42000,0,0,0,0,0,84000,0,0,0,0,0,-42000,0,0,0),
    @(0,-5040,0,8400,0,0,0,0,0,33600,0,0,0,5040,0,8400,0,0),
    @(10080,0,0,0,16800,0,0,0,0,0,67200,0,-10080,0,0,0,16800,0),
    @(0,0,0,0,0,-800,0,0,0,0,0,1600,0,0,0,0,0,-800),
    @(0,0,0,0,0,0,-4032,0,0,0,-10080,0,4032,0,0,0,-10080,0),
    @(0,0,0,0,0,0,0,-2016,0,5040,0,0,0,2016,0,5040,0,0),
    @(0,0,0,0,0,0,0,0,-42000,0,0,0,0,0,42000,0,0,0),
    @(0,0,0,0,0,0,0,-5040,0,8400,0,0,0,5040,0,16800,0,0),
    @(0,0,0,0,0,0,10080,0,0,0,16800,0,-10080,0,0,0,33600,0),
    @(0,0,0,0,0,0,0,0,0,0,0,-800,0,0,0,0,0,800)
  )

$wsBarra1 = $wb.Worksheets.Item("Barra 1 (K1)")
for ($r = 0; $r -lt $k12x12.Length; $r++) {
  $row = $k12x12[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $wsBarra1.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}

$wsBarra2 = $wb.Worksheets.Item("Barra 2 (K2)")
for ($r = 0; $r -lt $k12x12.Length; $r++) {
  $row = $k12x12[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $wsBarra2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}

$wsGlobal = $wb.Worksheets.Item("Global")
for ($r = 0; $r -lt $kGlobal18x18.Length; $r++) {
  $row = $kGlobal18x18[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $wsGlobal.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}

# ---------------------------------------------------------------------
# 2) Add the new sheets: F (load vector), Desplazamientos (nodal
#    displacements) and Reacciones (support reactions), right after
#    the existing "Global" sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsF = $wb.Worksheets.Add($null, $lastSheet)
$wsF.Name = "F"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDesplazamientos = $wb.Worksheets.Add($null, $lastSheet)
$wsDesplazamientos.Name = "Desplazamientos"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsReacciones = $wb.Worksheets.Add($null, $lastSheet)
$wsReacciones.Name = "Reacciones"

$colF = @(
    [double]"-3.673940397442059e-14",
    -648,
    [double]"-3.967855629237424e-14",
    0,
    [double]"-4.408728476930472e-14",
    -720,
    [double]"-6.123233995736766e-14",
    -1000,
    [double]"-6.123233995736768e-14",
    0,
    [double]"-1.469576158976824e-14",
    [double]"-239.9999999999999",
    [double]"-2.449293598294707e-14",
    [double]"-352.0000000000001",
    [double]"-2.155378366499342e-14",
    0,
    [double]"2.939152317953648e-14",
    [double]"480.0000000000001"
  )

$colDesplazamientos = @(
    0,
    0,
    0,
    0,
    0,
    0,
    [double]"-7.593296125665633e-18",
    [double]"-0.2480158730158731",
    [double]"-7.289564280639009e-19",
    0,
    [double]"-2.186869284191702e-19",
    [double]"-0.1499999999999999",
    0,
    0,
    0,
    0,
    0,
    0
  )

$colReacciones = @(
    [double]"6.515120971463919e-14",
    1148,
    [double]"7.029472627105808e-14",
    -1250,
    [double]"1.169537693185722e-13",
    840,
    0,
    0,
    0,
    0,
    0,
    0,
    [double]"5.731347020009613e-14",
    [double]"852.0000000000002",
    [double]"5.216995364367726e-14",
    1250,
    [double]"-1.096058885236881e-13",
    [double]"-360.0000000000002"
  )

for ($i = 0; $i -lt $colF.Length; $i++) {
  $wsF.Cells.Item($i + 1, 1).Value = $colF[$i]
}

for ($i = 0; $i -lt $colDesplazamientos.Length; $i++) {
  $wsDesplazamientos.Cells.Item($i + 1, 1).Value = $colDesplazamientos[$i]
}

for ($i = 0; $i -lt $colReacciones.Length; $i++) {
  $wsReacciones.Cells.Item($i + 1, 1).Value = $colReacciones[$i]
}
